# bias_tests.xlsx — "Further edits, table updates -Discarding old .tex files."
#
# Summary of the edit being applied:
#  - Four data rows (8, 9, 12, 13) got corrected Beta_Egger / SE(Beta_Egger) /
#    p_Egger / p_p-uniform numbers.
#  - All the numeric columns (D, E, F, G) are now right-aligned; column F
#    (p_Egger) additionally uses a ".000" number format for the numeric rows.
#  - The old p_{TES} column (H) was removed entirely.
#  - The selection cursor moved to E13.
#  - Page setup now explicitly records portrait orientation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Correct the four rows whose stats changed.
# ---------------------------------------------------------------------------
$ws.Range("D8").Value = 0.183
$ws.Range("E8").Value = 1.17
$ws.Range("F8").Value = 0.876
$ws.Range("G8").Value = 0.322

$ws.Range("D9").Value = 0.843
$ws.Range("E9").Value = 0.758
$ws.Range("F9").Value = 0.266
$ws.Range("G9").Value = 0.802

$ws.Range("D12").Value = -0.447
$ws.Range("E12").Value = 1.469
$ws.Range("F12").Value = 0.761
$ws.Range("G12").Value = 0.198

$ws.Range("D13").Value = 0.334
$ws.Range("E13").Value = 1.366
$ws.Range("F13").Value = 0.807
$ws.Range("G13").Value = 0.201

# ---------------------------------------------------------------------------
# 2. Right-align the numeric columns D:G for every data row (2-15).
#    Column F (p_Egger) also gets a ".000" number format everywhere it holds
#    a real number (rows 5-15; rows 2-3 hold the "< .001" text, row 4 holds
#    the "-" placeholder text).
# ---------------------------------------------------------------------------
$ws.Range("D2:G15").HorizontalAlignment = -4152   # xlRight

$ws.Range("F5:F15").NumberFormat = ".000"

# ---------------------------------------------------------------------------
# 3. Drop the old p_{TES} column entirely (column H).
# ---------------------------------------------------------------------------
$ws.Columns("H").Delete()

# ---------------------------------------------------------------------------
# 4. Cosmetic bits: selection cursor + page orientation.
# ---------------------------------------------------------------------------
$ws.Range("E13").Select()
$ws.PageSetup.Orientation = 1   # xlPortrait

Write-Output "edit applied"
